$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (style s="1") onto the new header cells
# AD1:AF1 so "Wins"/"Losses"/"Ties" match the existing header formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (Wins/Losses/Ties) is the same for every player row (2-55).
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 78
    $ws.Cells.Item($r, 32).Value = 0
}

"done"
